# Kaltura videos import - rewrite the data sheet to pull entries from the
# Kaltura server instead of the previous static video/audio sample rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the values of the previously used range but keep the existing cell
# styles (header shading/border, numeric date format on A:B) intact.
[void]$ws.Range("A1:I4").ClearContents()

# Columns H (Miniatura) and I (Tipo_de_archivo) are no longer produced by
# the Kaltura import, so drop them (value + style) entirely.
[void]$ws.Range("H1:I4").Clear()

# ---- Header row -----------------------------------------------------
$ws.Range("A1").Value = "Timestamp_inicio"
$ws.Range("B1").Value = "Timestamp_final"
$ws.Range("C1").Value = "Titulo"
$ws.Range("D1").Value = "Subtitulo"
$ws.Range("E1").Value = "Estudio"
$ws.Range("F1").Value = "Duracion"
$ws.Range("G1").Value = "Nombre_del_archivo"

# ---- Data rows --------------------------------------------------------
# Row 2 : Podcast de prueba de video
$ws.Range("A2").Value = 1557853200000
$ws.Range("B2").Value = 1557856800000
$ws.Range("C2").Value = "Podcast de prueba de vídeo"
$ws.Range("D2").Value = "Podcast de vídeo"
$ws.Range("E2").Value = "Campus Sur Radio"
$ws.Range("F2").Value = 127
$ws.Range("G2").Value = "0_c58wk4db"

# Row 3 : Ejemplo de video 2
$ws.Range("A3").Value = 1558023276000
$ws.Range("B3").Value = 1558023276000
$ws.Range("C3").Value = "Ejemplo de vídeo 2"
$ws.Range("D3").Value = "Subtítulo 2"
$ws.Range("E3").Value = "Campus Sur Radio"
$ws.Range("F3").Value = 127
$ws.Range("G3").Value = "0_8ba2t465"

# Row 4 : Ejemplo de video 3
$ws.Range("A4").Value = 1558023276000
$ws.Range("B4").Value = 1558023276000
$ws.Range("C4").Value = "Ejemplo de vídeo 3"
$ws.Range("D4").Value = "Subtítulo 3"
$ws.Range("E4").Value = "Campus Sur Radio"
$ws.Range("F4").Value = 128
$ws.Range("G4").Value = "0_sr0gkznr"

# Row 5 : Ejemplo de video 4
$ws.Range("A5").Value = 1558023276000
$ws.Range("B5").Value = 1558023276000
$ws.Range("C5").Value = "Ejemplo de vídeo 4"
$ws.Range("D5").Value = "Subtítulo 4"
$ws.Range("E5").Value = "Campus Sur Radio"
$ws.Range("F5").Value = 129
$ws.Range("G5").Value = "0_hvo9z7lv"

# Row 6 : Ejemplo de video 5
$ws.Range("A6").Value = 1558023276000
$ws.Range("B6").Value = 1558023276000
$ws.Range("C6").Value = "Ejemplo de vídeo 5"
$ws.Range("D6").Value = "Subtítulo 5"
$ws.Range("E6").Value = "Campus Sur Radio"
$ws.Range("F6").Value = 130
$ws.Range("G6").Value = "0_5r02f10i"

# Row 7 : Audio
$ws.Range("A7").Value = 1558091306400
$ws.Range("B7").Value = 1558089866400
$ws.Range("C7").Value = "Audio"
$ws.Range("D7").Value = "Esto es un audio"
$ws.Range("E7").Value = "Campus Sur Radio"
$ws.Range("F7").Value = 250
$ws.Range("G7").Value = "0_97s649gx"

# ---- Selection matches the saved state in the workbook ---------------
[void]$ws.Range("G3").Select()
